$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title: October -> November
$ws.Range("A1").Value = "Table 6.2.A. Net Summer Capacity of Utility Scale Units by Technology and by State, November 2016 and 2015 (Megawatts)"

# Row 3 period headers: October 2016/2015 -> November 2016/2015
# (NumberFormat is toggled to text then restored so Excel does not
#  auto-convert the "Month YYYY" string into a date serial, and so the
#  cell keeps its original style index.)
foreach ($col in @("B","D","F","H","J","L","N")) {
    $ws.Range($col + "3").NumberFormat = "@"
    $ws.Range($col + "3").Value = "November 2016"
    $ws.Range($col + "3").NumberFormat = "#,##0.0"
}
foreach ($col in @("C","E","G","I","K","M","O")) {
    $ws.Range($col + "3").NumberFormat = "@"
    $ws.Range($col + "3").Value = "November 2015"
    $ws.Range($col + "3").NumberFormat = "#,##0.0"
}

# Updated capacity figures (November 2016 run supersedes October 2016 run)
$ws.Range("C4").Value = 4973
$ws.Range("D4").Value = 22730.1
$ws.Range("E4").Value = 22753.3
$ws.Range("N4").Value = 33685.3
$ws.Range("O4").Value = 33569.7
$ws.Range("E5").Value = 6309.6
$ws.Range("O5").Value = 8784.2
$ws.Range("C7").Value = 950.3
$ws.Range("E7").Value = 9839.3
$ws.Range("O7").Value = 13219.9
$ws.Range("C9").Value = 57.3
$ws.Range("D9").Value = 1809.3
$ws.Range("N9").Value = 1878.6
$ws.Range("O9").Value = 1848.6
$ws.Range("C11").Value = 10644.3
$ws.Range("D11").Value = 69323.9
$ws.Range("E11").Value = 68015.6
$ws.Range("N11").Value = 102700.6
$ws.Range("O11").Value = 101344.7
$ws.Range("E12").Value = 13535.5
$ws.Range("O12").Value = 18757.2
$ws.Range("D13").Value = 25983
$ws.Range("E13").Value = 26310.5
$ws.Range("N13").Value = 39916.8
$ws.Range("O13").Value = 40243.3
$ws.Range("C14").Value = 2852.6
$ws.Range("O14").Value = 42344.2
$ws.Range("B15").Value = 10239.7
$ws.Range("C15").Value = 10049.9
$ws.Range("D15").Value = 113805.9
$ws.Range("I15").Value = 100.6
$ws.Range("N15").Value = 145226.4
$ws.Range("O15").Value = 147716.2
$ws.Range("C17").Value = 1992.5
$ws.Range("O17").Value = 26344.9
$ws.Range("B18").Value = 2257.7
$ws.Range("N18").Value = 29145.9
$ws.Range("I19").Value = 28
$ws.Range("O19").Value = 28742.6
$ws.Range("C20").Value = 1100.5
$ws.Range("D20").Value = 14340.1
$ws.Range("N20").Value = 16666.9
$ws.Range("O20").Value = 16925.3
$ws.Range("B21").Value = 22028.2
$ws.Range("C21").Value = 20527.3
$ws.Range("E21").Value = 61263.3
$ws.Range("N21").Value = 88710.5
$ws.Range("O21").Value = 88329.6
$ws.Range("B22").Value = 6467.3
$ws.Range("C22").Value = 6238.9
$ws.Range("N22").Value = 16635.9
$ws.Range("O22").Value = 16790.4
$ws.Range("E23").Value = 9583.1
$ws.Range("O23").Value = 13950.5
$ws.Range("E24").Value = 10240.6
$ws.Range("O24").Value = 15820.7
$ws.Range("B27").Value = 2891.5
$ws.Range("N27").Value = 7512.3
$ws.Range("B29").Value = 14993
$ws.Range("C29").Value = 13403.7
$ws.Range("E29").Value = 158279.9
$ws.Range("N29").Value = 207823.1
$ws.Range("O29").Value = 204742.6
$ws.Range("E32").Value = 53374.6
$ws.Range("O32").Value = 58719.2
$ws.Range("C33").Value = 3034.5
$ws.Range("O33").Value = 36148.6
$ws.Range("C34").Value = 1002.4
$ws.Range("O34").Value = 12391.9
$ws.Range("B35").Value = 4295.2
$ws.Range("C35").Value = 3454.6
$ws.Range("N35").Value = 31677.1
$ws.Range("O35").Value = 30836.5
$ws.Range("E36").Value = 11635.2
$ws.Range("O36").Value = 22697.8
$ws.Range("B37").Value = 1841.4
$ws.Range("E37").Value = 16845.8
$ws.Range("N37").Value = 26598.9
$ws.Range("O37").Value = 25409.8
$ws.Range("E38").Value = 14163.3
$ws.Range("O38").Value = 15114.8
$ws.Range("C39").Value = 8011.4
$ws.Range("E39").Value = 68051.3
$ws.Range("O39").Value = 87548.5
$ws.Range("C41").Value = 906
$ws.Range("E41").Value = 19153.5
$ws.Range("O41").Value = 20059.5
$ws.Range("B44").Value = 29828.1
$ws.Range("C44").Value = 25855.5
$ws.Range("N44").Value = 185186.5
$ws.Range("O44").Value = 181522.4
$ws.Range("C47").Value = 5368.7
$ws.Range("O47").Value = 24361.1
$ws.Range("B48").Value = 21160.9
$ws.Range("C48").Value = 18221.1
$ws.Range("N48").Value = 119306.6
$ws.Range("O48").Value = 116231.3
$ws.Range("B49").Value = 23818.6
$ws.Range("C49").Value = 21645.6
$ws.Range("N49").Value = 92744.6
$ws.Range("O49").Value = 90563.6
$ws.Range("B50").Value = 4652.6
$ws.Range("N50").Value = 28188.3
$ws.Range("B51").Value = 4095.9
$ws.Range("C51").Value = 3823
$ws.Range("N51").Value = 16030.7
$ws.Range("O51").Value = 15817.8
$ws.Range("C54").Value = 2554.4
$ws.Range("O54").Value = 10813.1
$ws.Range("B56").Value = 1478.6
$ws.Range("N56").Value = 8981.2
$ws.Range("B58").Value = 65092.3
$ws.Range("C58").Value = 62627.8
$ws.Range("D58").Value = 52496.4
$ws.Range("E58").Value = 51992.9
$ws.Range("N58").Value = 125331.8
$ws.Range("O58").Value = 122322.3
$ws.Range("B59").Value = 28203.2
$ws.Range("C59").Value = 25860.3
$ws.Range("D59").Value = 43308.9
$ws.Range("E59").Value = 43325.4
$ws.Range("N59").Value = 77774.2
$ws.Range("O59").Value = 75409.3
$ws.Range("B60").Value = 12102
$ws.Range("C60").Value = 12045.7
$ws.Range("N60").Value = 16479.8
$ws.Range("O60").Value = 15910.5
$ws.Range("B61").Value = 24787.1
$ws.Range("N61").Value = 31077.8
$ws.Range("D62").Value = 4182.2
$ws.Range("N62").Value = 5356.6
$ws.Range("D63").Value = 2125
$ws.Range("N63").Value = 2667.2
$ws.Range("B65").Value = 191164
$ws.Range("C65").Value = 178810.1
$ws.Range("D65").Value = 758464.1
$ws.Range("E65").Value = 761097.9
$ws.Range("I65").Value = 320.7
$ws.Range("N65").Value = 1073510.1
$ws.Range("O65").Value = 1062937.2
